# Locate the block of 14 trailing blank paragraphs (with pBdr borders) that sit
# at the very end of the document, right after the last
# "Drive to KAID and unload" paragraph. That's the block the diff rewrites:
#   - 7 of the "top+bottom border" blank paragraphs collapse into 1 paragraph
#     (keeps its pPr) that receives "LOCATION" / "NAME" text
#   - the 6 "bottom border only" blank paragraphs receive people/location text
#   - 2 brand-new "bottom border only" paragraphs are appended at the end
#     (Sheldon Hart, and a last mostly-empty one)

$d = $word.ActiveDocument

$TAB = "$([char]9)"

$n = $d.Paragraphs.Count
$anchorIdx = -1
for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Drive to KAID and unload*") {
        $anchorIdx = $i
    }
}

if ($anchorIdx -eq -1) {
    throw "Could not locate anchor paragraph 'Drive to KAID and unload'"
}

$firstBlankIdx = $anchorIdx + 1

# --- Delete 7 of the 8 "top+bottom border" blank paragraphs -----------------
# Keep the first one (it keeps its pPr/border and becomes the LOCATION/NAME
# row); remove the following 7 entirely.
for ($k = 0; $k -lt 7; $k++) {
    $p = $d.Paragraphs.Item($firstBlankIdx + 1)
    $p.Range.Delete()
}

# --- Append 2 brand-new paragraphs at the end --------------------------------
# They inherit the pPr (bottom-border-only tab row) of the current last
# paragraph, matching the target's new paragraphs.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()

# --- Fill in the 9 surviving/added paragraphs with their text ---------------
# Paragraph 1: LOCATION / NAME header row (keeps its top+bottom border pPr)
$p = $d.Paragraphs.Item($firstBlankIdx)
$p.Range.InsertBefore($TAB + "LOCATION" + $TAB + "NAME")

# Paragraph 2: ISHS row
$p = $d.Paragraphs.Item($firstBlankIdx + 1)
$p.Range.InsertBefore($TAB + "ISHS" + $TAB)

# Paragraph 3: POCATELLO / Shawn Phelps
$p = $d.Paragraphs.Item($firstBlankIdx + 2)
$p.Range.InsertBefore($TAB + "POCATELLO" + $TAB + "Shawn Phelps")

# Paragraph 4: FRANKLIN / Susan Hawkes
$p = $d.Paragraphs.Item($firstBlankIdx + 3)
$p.Range.InsertBefore($TAB + "FRANKLIN" + $TAB + "Susan Hawkes")

# Paragraph 5: (blank location) / Shawnee Hawkes
$p = $d.Paragraphs.Item($firstBlankIdx + 4)
$p.Range.InsertBefore($TAB + $TAB + "Shawnee Hawkes")

# Paragraph 6: RIGBY / Leon Guyman
$p = $d.Paragraphs.Item($firstBlankIdx + 5)
$p.Range.InsertBefore($TAB + "RIGBY" + $TAB + "Leon Guyman")

# Paragraph 7: (blank location) / Pat Lyn Scott
$p = $d.Paragraphs.Item($firstBlankIdx + 6)
$p.Range.InsertBefore($TAB + $TAB + "Pat Lyn Scott")

# Paragraph 8 (new): (blank location) / Sheldon Hart
$p = $d.Paragraphs.Item($firstBlankIdx + 7)
$p.Range.InsertBefore($TAB + $TAB + "Sheldon Hart")

# Paragraph 9 (new): single trailing tab, no other text
$p = $d.Paragraphs.Item($firstBlankIdx + 8)
$p.Range.InsertBefore($TAB)

Write-Host "Done. Total paragraphs now:" $d.Paragraphs.Count
